$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" date column (C) for existing data rows 2-23
#    from 2023-09-13 (45182) to 2023-09-15 (45184).
$lastRow = 23
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 23 regains the explicit 15pt row height that the other data rows have.
$ws.Rows.Item(23).RowHeight = 15

# 2. Append four new data rows (24-27) with the new avverkningsanmälningar.
$newRows = @(
    @{ Row = 24; A = "A 42957-2023"; G = 0.5 },
    @{ Row = 25; A = "A 42951-2023"; G = 0.8 },
    @{ Row = 26; A = "A 42960-2023"; G = 0.4 },
    @{ Row = 27; A = "A 42955-2023"; G = 0.3 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.A            # A - Beteckning
    $ws.Cells.Item($r, 2).Value = 45182              # B - Datum
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 3).Value = 45184              # C - Förändrad
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 4).Value = "VÄRMLANDS LÄN"    # D - Län
    $ws.Cells.Item($r, 5).Value = "HAMMARÖ"          # E - Kommun
    $ws.Cells.Item($r, 6).Value = "Övriga Aktiebolag" # F - Markägare
    $ws.Cells.Item($r, 7).Value = $item.G            # G - Area (ha)
    $ws.Cells.Item($r, 8).Value = 0                  # H - Fridlysta
    $ws.Cells.Item($r, 9).Value = 0                  # I - Signalarter
    $ws.Cells.Item($r, 10).Value = 0                 # J - NT
    $ws.Cells.Item($r, 11).Value = 0                 # K - VU
    $ws.Cells.Item($r, 12).Value = 0                 # L - EN
    $ws.Cells.Item($r, 13).Value = 0                 # M - CR
    $ws.Cells.Item($r, 14).Value = 0                 # N - RE
    $ws.Cells.Item($r, 15).Value = 0                 # O - Rödlistade
    $ws.Cells.Item($r, 16).Value = 0                 # P - Hotade
    $ws.Cells.Item($r, 17).Value = 0                 # Q - Alla arter
    $ws.Cells.Item($r, 18).WrapText = $true           # R - Artnamn (blank, wrapped)
}

# New rows 24-26 get the explicit 15pt row height (row 27, the new last row,
# is left without an explicit height, matching how a freshly appended last
# row behaves in this workbook).
$ws.Rows.Item(24).RowHeight = 15
$ws.Rows.Item(25).RowHeight = 15
$ws.Rows.Item(26).RowHeight = 15
